# ermia-arch.pptx — "paper: correct arch diagram"
#
# Moves/resizes the "Central Log buffer" box (now split onto two lines,
# "Centralized" / "log buffer") together with the "Commit" box and the
# right-arrow connector that points at it, nudging the whole cluster up a
# little and growing the log-buffer box to fit the extra line of text.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- "Central Log buffer" -> two centered lines: "Centralized" / "log buffer"
$rectLogBuffer = $s.Shapes.Item("Rectangle 20")
$rectLogBuffer.TextFrame.TextRange.Text = "Centralized`rlog buffer"

# Reposition/resize the box (EMU 220828,2902411 / 1670116x509404 expressed
# in points, the unit Shape.Left/Top/Width/Height use).
$rectLogBuffer.Left   = 17.388031996062992
$rectLogBuffer.Top    = 228.53629921259844
$rectLogBuffer.Width  = 131.5051968503937
$rectLogBuffer.Height = 40.11055188110236

# --- "Commit " box: shifts up slightly (no size/text change)
$rectCommit = $s.Shapes.Item("Rectangle 147")
$rectCommit.Left = 53.97338682677165
$rectCommit.Top  = 181.52244094488188

# --- Right-arrow connector feeding the log buffer: shifts up slightly too
$arrow = $s.Shapes.Item("Right Arrow 146")
$arrow.Left = 30.37283464566929
$arrow.Top  = 187.57827001653544
